$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.071.19"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "2.262.09"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'232.75"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'0.636"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("D7").Value = "'63.69"
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.449"
$ws.Range("E9").Value = "  +5.53%  "
$ws.Range("D10").Value = "'0.0980"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").Value = "'58.99"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").Value = "'26.48"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D14").Value = "2.593.26"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "'15.62"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("D16").Value = "'6.09"
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").Value = "'0.837"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").Value = "2.248.58"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").Value = "43.844.01"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "0.0₃0988"
$ws.Range("E20").Value = "  +4.08%  "
$ws.Range("D21").Value = "'73.57"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("D23").Value = "'250.14"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  -5.14%  "
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "'3.35"
$ws.Range("E26").Value = "  +21.82%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'2.22"
$ws.Range("E27").Value = "  -4.77%  "
$ws.Range("D28").Value = "'9.88"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "'173.36"
$ws.Range("D30").Value = "'21.77"
$ws.Range("E30").Value = "  +5.65%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "'1.44"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("D34").Value = "'4.86"
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("D35").Value = "'0.0685"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("E36").Value = "  -5.18%  "
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("D38").Value = "'6.43"
$ws.Range("E38").Value = "  -5.03%  "
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'8.66"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("D43").Value = "'0.000223"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'17.20"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'98.14"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("E46").Value = "  -2.66%  "
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").Value = "'4.35"
$ws.Range("E48").Value = "  -8.18%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.34"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.446.47"
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'1.09"
$ws.Range("E51").Value = "  -0.29%  "
